$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 1
    3 = 3
    5 = 2
    6 = 3
    7 = 2
    8 = 2
    9 = 2
    11 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 3
    16 = 3
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 2
    23 = 1
    24 = 3
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 3
    36 = 2
    37 = 1
    38 = 2
    39 = 1
    40 = 1
    43 = 3
    44 = 3
    45 = 1
    46 = 3
    47 = 1
    48 = 2
    49 = 2
    50 = 2
    51 = 1
    52 = 1
    54 = 2
    55 = 2
    56 = 2
    57 = 1
    58 = 1
    59 = 3
    60 = 2
    61 = 2
    62 = 1
    63 = 1
    64 = 1
    65 = 1
    66 = 1
    67 = 2
    68 = 2
    69 = 2
    70 = 1
    71 = 1
    73 = 1
    74 = 3
    78 = 1
    79 = 1
    80 = 1
    81 = 2
    91 = 2
    92 = 1
    93 = 1
    94 = 1
    96 = 1
    97 = 1
    99 = 1
    101 = 1
    108 = 1
    109 = 1
    113 = 1
    114 = 1
    115 = 1
    116 = 1
    117 = 1
    118 = 1
    121 = 1
    127 = 1
    130 = 1
    132 = 1
    134 = 1
    136 = 1
    137 = 1
    138 = 1
    139 = 1
    140 = 2
    146 = 1
    147 = 1
    150 = 1
    153 = 1
    154 = 1
    155 = 1
    157 = 1
    159 = 1
    160 = 1
    161 = 1
    162 = 1
    163 = 1
    164 = 1
    165 = 1
    166 = 1
    167 = 1
    168 = 1
    169 = 1
    171 = 1
    172 = 1
    173 = 1
    177 = 1
    179 = 1
    183 = 1
    184 = 1
    185 = 1
    190 = 1
    191 = 1
    194 = 1
    195 = 1
    200 = 1
    201 = 1
    203 = 1
    206 = 1
    209 = 1
    210 = 1
    212 = 1
    213 = 1
    215 = 1
    216 = 1
    218 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}
